$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay text even when the new value looks like a
    # plain number (e.g. "1.44"), then restore the original (default)
    # cell style so no stray formatting is introduced.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "97.770.25"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "3.378.32"
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue $ws.Range("D5") "252.01"
$ws.Range("E5").Value = "  -1.90%  "

Set-TextValue $ws.Range("D6") "668.58"
$ws.Range("E6").Value = "  +0.68%  "

Set-TextValue $ws.Range("D7") "1.44"
$ws.Range("E7").Value = "  -6.34%  "

Set-TextValue $ws.Range("D8") "0.426"
$ws.Range("E8").Value = "  -10.16%  "

$ws.Range("E9").Value = "  -0.04%  "

Set-TextValue $ws.Range("D10") "1.04"
$ws.Range("E10").Value = "  -3.59%  "

$ws.Range("D11").Value = "3.376.61"
$ws.Range("E11").Value = "  +0.86%  "

$ws.Range("E12").Value = "  -1.67%  "

Set-TextValue $ws.Range("D13") "41.47"
$ws.Range("E13").Value = "  -1.56%  "

$ws.Range("D14").Value = "97.499.87"
$ws.Range("E14").Value = "  -0.52%  "

Set-TextValue $ws.Range("D15") "6.20"
$ws.Range("E15").Value = "  +8.95%  "

Set-TextValue $ws.Range("D16") "0.0000260"
$ws.Range("E16").Value = "  -5.39%  "

$ws.Range("D17").Value = "4.005.50"
$ws.Range("E17").Value = "  +0.83%  "

Set-TextValue $ws.Range("D18") "8.72"
$ws.Range("E18").Value = "  +14.02%  "

$ws.Range("D19").Value = "3.370.83"
$ws.Range("E19").Value = "  +0.20%  "

Set-TextValue $ws.Range("D20") "0.575"
$ws.Range("E20").Value = "  +32.48%  "

Set-TextValue $ws.Range("D21") "17.17"
$ws.Range("E21").Value = "  +2.64%  "

Set-TextValue $ws.Range("D22") "11.01"
$ws.Range("E22").Value = "  +4.29%  "

Set-TextValue $ws.Range("D23") "508.05"
$ws.Range("E23").Value = "  -3.87%  "

Set-TextValue $ws.Range("D24") "3.41"
$ws.Range("E24").Value = "  -4.72%  "

Set-TextValue $ws.Range("D25") "0.0000203"
$ws.Range("E25").Value = "  -6.79%  "

Set-TextValue $ws.Range("D26") "6.45"
$ws.Range("E26").Value = "  +4.35%  "

Set-TextValue $ws.Range("D27") "99.25"
$ws.Range("E27").Value = "  -3.05%  "

Set-TextValue $ws.Range("D28") "12.44"
$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("D29").Value = "3.567.20"
$ws.Range("E29").Value = "  +1.02%  "

Set-TextValue $ws.Range("D30") "0.151"
$ws.Range("E30").Value = "  +1.90%  "

Set-TextValue $ws.Range("D31") "11.37"
$ws.Range("E31").Value = "  +3.21%  "

Set-TextValue $ws.Range("D32") "0.997"
$ws.Range("E32").Value = "  -0.15%  "

Set-TextValue $ws.Range("D33") "0.192"
$ws.Range("E33").Value = "  +1.32%  "

Set-TextValue $ws.Range("D34") "2.61"
$ws.Range("E34").Value = "  +23.70%  "

Set-TextValue $ws.Range("D35") "1.00"
$ws.Range("E35").Value = "  +0.28%  "

Set-TextValue $ws.Range("D36") "0.563"
$ws.Range("E36").Value = "  +5.05%  "

Set-TextValue $ws.Range("D37") "29.08"
$ws.Range("E37").Value = "  -0.89%  "

Set-TextValue $ws.Range("D38") "7.93"
$ws.Range("E38").Value = "  +1.71%  "

Set-TextValue $ws.Range("D39") "1.50"
$ws.Range("E39").Value = "  +12.54%  "

Set-TextValue $ws.Range("D40") "536.48"
$ws.Range("E40").Value = "  +2.41%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("E42").Value = "  -3.81%  "

Set-TextValue $ws.Range("D43") "24.70"
$ws.Range("E43").Value = "  -0.01%  "

Set-TextValue $ws.Range("D44") "9.09"
$ws.Range("E44").Value = "  +16.94%  "

Set-TextValue $ws.Range("D45") "0.856"
$ws.Range("E45").Value = "  +4.42%  "

Set-TextValue $ws.Range("D46") "0.0431"
$ws.Range("E46").Value = "  -0.72%  "

Set-TextValue $ws.Range("D47") "3.70"
$ws.Range("E47").Value = "  -4.85%  "

Set-TextValue $ws.Range("D48") "5.71"
$ws.Range("E48").Value = "  +11.76%  "

Set-TextValue $ws.Range("D49") "1.71"
$ws.Range("E49").Value = "  +11.49%  "

Set-TextValue $ws.Range("D50") "54.34"
$ws.Range("E50").Value = "  +6.85%  "

Set-TextValue $ws.Range("D51") "3.20"
$ws.Range("E51").Value = "  -6.23%  "
